$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 value to the new predicted hardness value
$ws.Range("B3").Value = 68.518310546875

# Remove the now-obsolete composition rows (4-6)
$ws.Range("A4:B6").ClearContents()
